$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Y" flag in A2 was removed (cell cleared), which also drops the
# now-unused "Y" shared string and shifts the selection to A2.
$ws.Range("A2").ClearContents()
$ws.Range("A2").Select()
